$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172892689704895
$ws.Range("B1").Value = 2.36356520652771
$ws.Range("C1").Value = 5.020647525787354
$ws.Range("D1").Value = 2.341699361801147
$ws.Range("E1").Value = 1.22899329662323
